# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / Leve profit columns (H:N) for 43 rows
# across the 8 job sheets, reflecting refreshed market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29: H29=538.8889, I29=43.75, J29=4500, K29=131.25, L29=13500, M29=149.75, N29=-14062
$ws.Range("H29").Value = 538.8889
$ws.Range("I29").Value = 43.75
$ws.Range("J29").Value = 4500
$ws.Range("K29").Value = 131.25
$ws.Range("L29").Value = 13500
$ws.Range("M29").Value = 149.75
$ws.Range("N29").Value = -14062

# Row 70: H70=37042.043, J70=75397.10000000001, L70=226191.3, N70=-226731.3
$ws.Range("H70").Value = 37042.043
$ws.Range("J70").Value = 75397.10000000001
$ws.Range("L70").Value = 226191.3
$ws.Range("N70").Value = -226731.3

# Row 73: H73=37042.043, J73=75397.10000000001, L73=226191.3, N73=-228063.3
$ws.Range("H73").Value = 37042.043
$ws.Range("J73").Value = 75397.10000000001
$ws.Range("L73").Value = 226191.3
$ws.Range("N73").Value = -228063.3

# Row 98: H98=1027.6, I98=982.7059, J98=1282, K98=982.7059, L98=1282, M98=515.2941, N98=-4278
$ws.Range("H98").Value = 1027.6
$ws.Range("I98").Value = 982.7059
$ws.Range("J98").Value = 1282
$ws.Range("K98").Value = 982.7059
$ws.Range("L98").Value = 1282
$ws.Range("M98").Value = 515.2941
$ws.Range("N98").Value = -4278

# Row 113: H113=4700, I113=4700, K113=4700, M113=-1446
$ws.Range("H113").Value = 4700
$ws.Range("I113").Value = 4700
$ws.Range("K113").Value = 4700
$ws.Range("M113").Value = -1446

# Row 122: H122=1027.6, I122=982.7059, J122=1282, K122=2948.1177, L122=3846, M122=-498.1177000000002, N122=-8746
$ws.Range("H122").Value = 1027.6
$ws.Range("I122").Value = 982.7059
$ws.Range("J122").Value = 1282
$ws.Range("K122").Value = 2948.1177
$ws.Range("L122").Value = 3846
$ws.Range("M122").Value = -498.1177000000002
$ws.Range("N122").Value = -8746

# Row 127: H127=3656.5833, I127=3656.5833, K127=10969.7499, M127=-6009.749899999999
$ws.Range("H127").Value = 3656.5833
$ws.Range("I127").Value = 3656.5833
$ws.Range("K127").Value = 10969.7499
$ws.Range("M127").Value = -6009.749899999999

# Row 132: H132=2752.3333, I132=1378.2778, K132=4134.8334, M132=-1604.8334
$ws.Range("H132").Value = 2752.3333
$ws.Range("I132").Value = 1378.2778
$ws.Range("K132").Value = 4134.8334
$ws.Range("M132").Value = -1604.8334

$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2=1720.4615, I2=1130.3334, J2=3048.25, K2=1130.3334, L2=3048.25, M2=-1017.3334, N2=-3274.25
$ws.Range("H2").Value = 1720.4615
$ws.Range("I2").Value = 1130.3334
$ws.Range("J2").Value = 3048.25
$ws.Range("K2").Value = 1130.3334
$ws.Range("L2").Value = 3048.25
$ws.Range("M2").Value = -1017.3334
$ws.Range("N2").Value = -3274.25

# Row 13: H13=6499, I13=0, J13=6499, K13=0, L13=6499, N13=-6787
$ws.Range("H13").Value = 6499
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 6499
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 6499
$ws.Range("N13").Value = -6787
$ws.Range("M13").ClearContents()

# Row 32: H32=3627.2856, I32=2225.6216, K32=2225.6216, M32=-1938.6216
$ws.Range("H32").Value = 3627.2856
$ws.Range("I32").Value = 2225.6216
$ws.Range("K32").Value = 2225.6216
$ws.Range("M32").Value = -1938.6216

# Row 45: H45=3111.7273, I45=2049.2, K45=2049.2, M45=-1672.2
$ws.Range("H45").Value = 3111.7273
$ws.Range("I45").Value = 2049.2
$ws.Range("K45").Value = 2049.2
$ws.Range("M45").Value = -1672.2

# Row 88: H88=1636, I88=870.75, K88=870.75, M88=-464.75
$ws.Range("H88").Value = 1636
$ws.Range("I88").Value = 870.75
$ws.Range("K88").Value = 870.75
$ws.Range("M88").Value = -464.75

# Row 91: H91=1636, I91=870.75, K91=870.75, M91=533.25
$ws.Range("H91").Value = 1636
$ws.Range("I91").Value = 870.75
$ws.Range("K91").Value = 870.75
$ws.Range("M91").Value = 533.25

# Row 116: H116=1720.4615, I116=1130.3334, J116=3048.25, K116=1130.3334, L116=3048.25, M116=1163.6666, N116=-7636.25
$ws.Range("H116").Value = 1720.4615
$ws.Range("I116").Value = 1130.3334
$ws.Range("J116").Value = 3048.25
$ws.Range("K116").Value = 1130.3334
$ws.Range("L116").Value = 3048.25
$ws.Range("M116").Value = 1163.6666
$ws.Range("N116").Value = -7636.25

# Row 132: H132=2147.3447, I132=2145.16, J132=2161, K132=6435.48, L132=6483, M132=-3905.48, N132=-11543
$ws.Range("H132").Value = 2147.3447
$ws.Range("I132").Value = 2145.16
$ws.Range("J132").Value = 2161
$ws.Range("K132").Value = 6435.48
$ws.Range("L132").Value = 6483
$ws.Range("M132").Value = -3905.48
$ws.Range("N132").Value = -11543

$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3=1720.4615, I3=1130.3334, J3=3048.25, K3=1130.3334, L3=3048.25, M3=-1016.3334, N3=-3276.25
$ws.Range("H3").Value = 1720.4615
$ws.Range("I3").Value = 1130.3334
$ws.Range("J3").Value = 3048.25
$ws.Range("K3").Value = 1130.3334
$ws.Range("L3").Value = 3048.25
$ws.Range("M3").Value = -1016.3334
$ws.Range("N3").Value = -3276.25

# Row 20: H20=2037.8667, I20=2537.2, K20=2537.2, M20=-2290.2
$ws.Range("H20").Value = 2037.8667
$ws.Range("I20").Value = 2537.2
$ws.Range("K20").Value = 2537.2
$ws.Range("M20").Value = -2290.2

# Row 50: H50=145666, J50=145666, L50=145666, N50=-146814
$ws.Range("H50").Value = 145666
$ws.Range("J50").Value = 145666
$ws.Range("L50").Value = 145666
$ws.Range("N50").Value = -146814

# Row 64: H64=1166.5, J64=1099.8, L64=1099.8, N64=-1549.8
$ws.Range("H64").Value = 1166.5
$ws.Range("J64").Value = 1099.8
$ws.Range("L64").Value = 1099.8
$ws.Range("N64").Value = -1549.8

# Row 67: H67=1166.5, J67=1099.8, L67=1099.8, N67=-2659.8
$ws.Range("H67").Value = 1166.5
$ws.Range("J67").Value = 1099.8
$ws.Range("L67").Value = 1099.8
$ws.Range("N67").Value = -2659.8

# Row 134: H134=1212.9546, I134=983.55, K134=2950.65, M134=-415.6499999999996
$ws.Range("H134").Value = 1212.9546
$ws.Range("I134").Value = 983.55
$ws.Range("K134").Value = 2950.65
$ws.Range("M134").Value = -415.6499999999996

$ws = $wb.Worksheets.Item("CRP")
# Row 22: H22=56536.047, I22=131111, K22=131111, M22=-130761
$ws.Range("H22").Value = 56536.047
$ws.Range("I22").Value = 131111
$ws.Range("K22").Value = 131111
$ws.Range("M22").Value = -130761

# Row 31: H31=3865.3333, I31=3956.5715, J31=3785.5, K31=3956.5715, L31=3785.5, M31=-3661.5715, N31=-4375.5
$ws.Range("H31").Value = 3865.3333
$ws.Range("I31").Value = 3956.5715
$ws.Range("J31").Value = 3785.5
$ws.Range("K31").Value = 3956.5715
$ws.Range("L31").Value = 3785.5
$ws.Range("M31").Value = -3661.5715
$ws.Range("N31").Value = -4375.5

# Row 34: H34=3865.3333, I34=3956.5715, J34=3785.5, K34=3956.5715, L34=3785.5, M34=-3754.5715, N34=-4189.5
$ws.Range("H34").Value = 3865.3333
$ws.Range("I34").Value = 3956.5715
$ws.Range("J34").Value = 3785.5
$ws.Range("K34").Value = 3956.5715
$ws.Range("L34").Value = 3785.5
$ws.Range("M34").Value = -3754.5715
$ws.Range("N34").Value = -4189.5

# Row 62: H62=47562.668, I62=3554, J62=102573.5, K62=3554, L62=102573.5, M62=-2930, N62=-103821.5
$ws.Range("H62").Value = 47562.668
$ws.Range("I62").Value = 3554
$ws.Range("J62").Value = 102573.5
$ws.Range("K62").Value = 3554
$ws.Range("L62").Value = 102573.5
$ws.Range("M62").Value = -2930
$ws.Range("N62").Value = -103821.5

# Row 65: H65=47562.668, I65=3554, J65=102573.5, K65=17770, L65=512867.5, M65=-14650, N65=-519107.5
$ws.Range("H65").Value = 47562.668
$ws.Range("I65").Value = 3554
$ws.Range("J65").Value = 102573.5
$ws.Range("K65").Value = 17770
$ws.Range("L65").Value = 512867.5
$ws.Range("M65").Value = -14650
$ws.Range("N65").Value = -519107.5

$ws = $wb.Worksheets.Item("CUL")
# Row 12: H12=117.045456, J12=114.416664, L12=343.249992, N12=-689.249992
$ws.Range("H12").Value = 117.045456
$ws.Range("J12").Value = 114.416664
$ws.Range("L12").Value = 343.249992
$ws.Range("N12").Value = -689.249992

# Row 34: H34=4211.5, J34=4384.5713, L34=13153.7139, N34=-13321.7139
$ws.Range("H34").Value = 4211.5
$ws.Range("J34").Value = 4384.5713
$ws.Range("L34").Value = 13153.7139
$ws.Range("N34").Value = -13321.7139

# Row 39: H39=9416.0625, J39=12715.7, L39=38147.10000000001, N39=-38735.10000000001
$ws.Range("H39").Value = 9416.0625
$ws.Range("J39").Value = 12715.7
$ws.Range("L39").Value = 38147.10000000001
$ws.Range("N39").Value = -38735.10000000001

# Row 55: H55=3531.125, J55=3357, L55=10071, N55=-10425
$ws.Range("H55").Value = 3531.125
$ws.Range("J55").Value = 3357
$ws.Range("L55").Value = 10071
$ws.Range("N55").Value = -10425

# Row 132: H132=2108.6, I132=2120.6667, J132=2000, K132=19086.0003, L132=18000, M132=-16556.0003, N132=-23060
$ws.Range("H132").Value = 2108.6
$ws.Range("I132").Value = 2120.6667
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 19086.0003
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -16556.0003
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
# Row 3: H3=161, I3=170.25, J3=50, K3=170.25, L3=50, M3=-54.25, N3=-282
$ws.Range("H3").Value = 161
$ws.Range("I3").Value = 170.25
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 170.25
$ws.Range("L3").Value = 50
$ws.Range("M3").Value = -54.25
$ws.Range("N3").Value = -282

# Row 26: H26=8000, J26=0, L26=0
$ws.Range("H26").Value = 8000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

# Row 50: H50=8000, J50=0, L50=0
$ws.Range("H50").Value = 8000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# Row 102: H102=998.875, I102=622.75, J102=1375, K102=622.75, L102=1375, M102=999.25, N102=-4619
$ws.Range("H102").Value = 998.875
$ws.Range("I102").Value = 622.75
$ws.Range("J102").Value = 1375
$ws.Range("K102").Value = 622.75
$ws.Range("L102").Value = 1375
$ws.Range("M102").Value = 999.25
$ws.Range("N102").Value = -4619

# Row 122: H122=77754.64, I122=5415.091, K122=16245.273, M122=-13795.273
$ws.Range("H122").Value = 77754.64
$ws.Range("I122").Value = 5415.091
$ws.Range("K122").Value = 16245.273
$ws.Range("M122").Value = -13795.273

# Row 132: H132=0, I132=0, J132=0, K132=0, L132=0
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 132: H132=252000.5, I132=252000.5, K132=756001.5, M132=-753471.5
$ws.Range("H132").Value = 252000.5
$ws.Range("I132").Value = 252000.5
$ws.Range("K132").Value = 756001.5
$ws.Range("M132").Value = -753471.5

$ws = $wb.Worksheets.Item("WVR")
# Row 8: H8=2749.5, I8=0, J8=2749.5, K8=0, L8=2749.5, N8=-3029.5
$ws.Range("H8").Value = 2749.5
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2749.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2749.5
$ws.Range("N8").Value = -3029.5
$ws.Range("M8").ClearContents()

# Row 113: H113=750, I113=750, J113=0, K113=2250, L113=0, M113=-80
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -80
$ws.Range("N113").ClearContents()

# Row 119: H119=39773.25, J119=39773.25, L119=39773.25, N119=-49449.25
$ws.Range("H119").Value = 39773.25
$ws.Range("J119").Value = 39773.25
$ws.Range("L119").Value = 39773.25
$ws.Range("N119").Value = -49449.25

# Row 132: H132=0, I132=0, J132=0, K132=0, L132=0
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
